$d = $word.ActiveDocument

# The notice currently reads "...Licence expired on March 31st." and needs
# to read "...Licence expires on March 31st." (present tense). Locate the
# word "expired" and change only its final letter ("d" -> "s") so the
# surrounding text/runs are disturbed as little as possible.
$findRng = $d.Content.Duplicate
$findRng.Find.Execute("expired", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($findRng.Find.Found) {
    $lastChar = $d.Range($findRng.End - 1, $findRng.End)
    $lastChar.Text = "s"
}
